$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q1" right after "2021-Q4" (and before
#    "总计"), mirroring the per-quarter fund-holdings sheets already present.
# ---------------------------------------------------------------------------
$prevSheet = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $prevSheet)
$ws.Name = "2022-Q1"

# --- header row ---------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Copy the header formatting (bold, centered, bordered) from a sibling
# quarter sheet so the new tab matches the rest of the workbook.
$prevSheet.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# --- fund-holding rows ----------------------------------------------------
# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "210009", "金鹰核心资源混合",   "3.86", "94.96", "4.62", "0.1783", 8),
    @(1, "001167", "金鹰科技创新股票",   "4.03", "94.55", "4.42", "0.1781", 8),
    @(2, "162102", "金鹰中小盘精选混合", "4.60", "76.52", "3.62", "0.1665", 7)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Columns B,D,E,F,G hold numeric-looking text in the source data; force
    # text formatting BEFORE assigning so Excel doesn't silently coerce them
    # to numbers (and so leading zeros in fund codes, e.g. 001167, survive).
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]

    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row[5]

    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row[6]

    # 仓位排名 stays a real number.
    $ws.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# Copy the index-column (A) formatting from a sibling quarter sheet onto the
# three new index cells (same bold/bordered style as the header).
$prevSheet.Range("A2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down by one and renumbering the index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D2").EntireRow.Insert()

# The row-insert can drag the header's formatting onto the new row; the
# data columns (B:D) in every other row are unstyled, so strip it back off.
$total.Range("B2:D2").ClearFormats()

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.52

# Copy the index-column (A) formatting down into the newly inserted row, and
# renumber the index column (0..5) for every row, now that one more row
# exists.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

Write-Output "2022-Q1 sheet added and 总计 updated"
